$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Variable" column (H) previously stored an Integrated-Assessment-Model
# style path ("Capacity Additions|Electricity|<Variable_simplified>"). That
# IAM-specific prefix is dropped so the column just repeats the
# "Variable_simplified" (column G) value.
$ws.Range("H20").Value = "Biomass|w/o CCS"
$ws.Range("H21").Value = "Biomass|w/ CCS"
$ws.Range("H22").Value = "Coal|w/o CCS"
$ws.Range("H23").Value = "Coal|w/ CCS"
$ws.Range("H24").Value = "Geothermal"
$ws.Range("H25").Value = "Oil|w/o CCS"
$ws.Range("H26").Value = "Oil|w/o CCS"
$ws.Range("H27").Value = "Gas|w/o CCS"
$ws.Range("H28").Value = "Gas|w/o CCS"
$ws.Range("H29").Value = "Gas|w/ CCS"
$ws.Range("H30").Value = "Solar"
$ws.Range("H31").Value = "Solar"
$ws.Range("H32").Value = "Solar"
$ws.Range("H34").Value = "Hydro"
$ws.Range("H36").Value = "Wind|Onshore"
$ws.Range("H37").Value = "Wind|Offshore"
$ws.Range("H38").Value = "Nuclear"
$ws.Range("H43").Value = "Oil|w/o CCS"
$ws.Range("H44").Value = "Solar"
$ws.Range("H45").Value = "Solar"
$ws.Range("H79").Value = "Solar"
$ws.Range("H80").Value = "Wind|Onshore"
$ws.Range("H81").Value = "Storage Capacity"
$ws.Range("H82").Value = "Oil|w/ CCS"

# F82 was a stray duplicate of E82 ("Oil+ CCS") - drop the cell entirely.
$ws.Range("F82").ClearContents()

# Drop the stray trailing placeholder row at the bottom of the sheet so the
# used range (and dimension) shrinks down to the real data (A1:J87).
$ws.Rows.Item(1048576).Delete()

# Restore the selection/view to where the edit ended up.
$ws.Range("F82").Select()
